# 自动更新Excel文件脚本
# 对每一行（剩余天数列 E），每天递减 1；
# 当剩余天数递减到 0（即当前值为 1）时，重置为 10，
# 并将开始时间列 F（yyyyMMdd 数字）向后推 10 天，开启新的周期。
# 若某行的开始时间格式异常（非标准 yyyyMMdd 8位数字），则整行跳过，不做任何修改。

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colRemain = 5   # E列：剩余
$colStart  = 6   # F列：开始时间

for ($row = 2; $row -le $lastRow; $row++) {
    $remainCell = $ws.Cells.Item($row, $colRemain)
    $remainValue = $remainCell.Value()

    if ($remainValue -eq $null -or $remainValue -eq "") {
        continue
    }

    $startCell = $ws.Cells.Item($row, $colStart)
    $startValue2 = $startCell.Value()

    if ($startValue2 -eq $null -or $startValue2 -eq "") {
        continue
    }

    $startValue = [string]([int]$startValue2)

    # 开始时间格式异常（非8位 yyyyMMdd）时跳过该行，不做任何修改
    if ($startValue.Length -ne 8) {
        continue
    }

    $remain = [int]$remainValue

    if ($remain -eq 1) {
        # 周期结束，重置剩余天数并顺延开始时间 10 天
        $year  = [int]$startValue.Substring(0, 4)
        $month = [int]$startValue.Substring(4, 2)
        $day   = [int]$startValue.Substring(6, 2)

        $startDate = Get-Date -Year $year -Month $month -Day $day
        $newDate = $startDate.AddDays(10)

        $newStartValue = [int]$newDate.ToString("yyyyMMdd")

        $remainCell.Value = 10
        $startCell.Value = $newStartValue
    }
    else {
        $remainCell.Value = $remain - 1
    }
}

Write-Host "Update complete"
